# Populate the remaining box-score columns (D,F,G,I,J,K,L,M) for the
# "Thu, Feb 8, 2024" slate (rows 150-158 on Sheet1), which previously only
# had the schedule columns (A,B,C,E,H) filled in. Column M extends the
# existing ABS(D-F) "Diff" shared formula down through row 158, and the
# dependent summary formulas on Sheet2 (COUNTIFS / AVERAGEIFS / the
# cumulative array formula) recalc automatically off of that.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$rows = @(
    [PSCustomObject]@{ Row=150; D=131; F=109; G="NA";  I="Golden State Warriors";   J="Indiana Pacers";        K="Indiana Pacers";        KWin=$false; L="No"  },
    [PSCustomObject]@{ Row=151; D=111; F=127; G="NA";  I="Orlando Magic";           J="San Antonio Spurs";     K="Orlando Magic";         KWin=$true;  L="Yes" },
    [PSCustomObject]@{ Row=152; D=118; F=95;  G="NA";  I="Cleveland Cavaliers";     J="Brooklyn Nets";         K="Cleveland Cavaliers";   KWin=$true;  L="Yes" },
    [PSCustomObject]@{ Row=153; D=122; F=108; G="NA";  I="Dallas Mavericks";        J="New York Knicks";       K="New York Knicks";       KWin=$false; L="No"  },
    [PSCustomObject]@{ Row=154; D=118; F=110; G="NA";  I="Chicago Bulls";           J="Memphis Grizzlies";     K="Chicago Bulls";         KWin=$true;  L="Yes" },
    [PSCustomObject]@{ Row=155; D=129; F=105; G="NA";  I="Minnesota Timberwolves";  J="Milwaukee Bucks";       K="Milwaukee Bucks";       KWin=$false; L="No"  },
    [PSCustomObject]@{ Row=156; D=115; F=129; G="NA";  I="Phoenix Suns";            J="Utah Jazz";             K="Utah Jazz";             KWin=$false; L="No"  },
    [PSCustomObject]@{ Row=157; D=114; F=106; G="NA";  I="Denver Nuggets";          J="Los Angeles Lakers";    K="Denver Nuggets";        KWin=$true;  L="Yes" },
    [PSCustomObject]@{ Row=158; D=128; F=122; G="Yes"; I="Detroit Pistons";         J="Portland Trail Blazers"; K="Portland Trail Blazers"; KWin=$false; L="No" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws1.Cells.Item($row, 4).Value = $r.D   # D: Away Pts
    $ws1.Cells.Item($row, 6).Value = $r.F   # F: Home Pts
    $ws1.Cells.Item($row, 7).Value = $r.G   # G: Overtime
    $ws1.Cells.Item($row, 9).Value = $r.I   # I: Win
    $ws1.Cells.Item($row, 10).Value = $r.J  # J: Loss
    $ws1.Cells.Item($row, 11).Value = $r.K  # K: Forecasted

    if ($r.KWin) {
        $ws1.Cells.Item($row, 11).Interior.Color = 5287936
    }

    $ws1.Cells.Item($row, 12).Value = $r.L  # L: Correct

    $mCell = $ws1.Cells.Item($row, 13)      # M: Diff
    $mCell.Formula = "=ABS(D$row-F$row)"
    $mCell.NumberFormat = "#,##0"
}

# The user had been working on Sheet2 (selecting D26) before switching back
# to Sheet1, which is now the active/selected tab.
$ws2 = $wb.Worksheets.Item("Sheet2")
[void]$ws2.Select()
[void]$ws2.Range("D26").Select()
[void]$ws1.Select()
